# Adds a new "2022-Q3" sheet (fund holdings detail) right after the summary
# sheet "总计" and before the existing "2022-Q2" sheet, and inserts a
# corresponding summary row at the top of the "总计" sheet's data, shifting
# the existing summary rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a value into a cell while forcing TEXT storage (matches
# source data where numeric-looking strings like "4.28" or zero-padded
# fund codes like "014016" must stay text, not get coerced to numbers).
# ---------------------------------------------------------------------
function Set-TextValue($cell, $val) {
    $cell.Value = "'" + $val
}

# ---------------------------------------------------------------------
# Helper: write a cell's value while cloning the style (number format /
# font / border / alignment) from a reference cell that already carries
# the desired style index. Copy first (brings style + old value), then
# overwrite with the real value so the style sticks but the data is ours.
# ---------------------------------------------------------------------
function Set-StyledValue($styleSourceCell, $destCell, $val) {
    $styleSourceCell.Copy($destCell)
    $destCell.Value = $val
}

function Set-StyledTextValue($styleSourceCell, $destCell, $val) {
    $styleSourceCell.Copy($destCell)
    Set-TextValue $destCell $val
}

# =======================================================================
# 1. "总计" sheet: insert a new row 2 ("2022-Q3") and push the old rows
#    (2022-Q2, 2022-Q1, 2021-Q4) down by one.
# =======================================================================
$summary = $wb.Worksheets.Item(1)

# Make room: shift existing rows 2..4 down to 3..5 (bottom-up so we never
# clobber a row before reading it). Column A keeps its bold/centered/
# bordered "index" style (style source = the row directly above, which
# already carries that style).
Set-StyledValue $summary.Cells.Item(4,1) $summary.Cells.Item(5,1) 3
$summary.Cells.Item(5,2).Value = $summary.Cells.Item(4,2).Value2
$summary.Cells.Item(5,3).Value = $summary.Cells.Item(4,3).Value2
$summary.Cells.Item(5,4).Value = $summary.Cells.Item(4,4).Value2

$summary.Cells.Item(4,2).Value = $summary.Cells.Item(3,2).Value2
$summary.Cells.Item(4,3).Value = $summary.Cells.Item(3,3).Value2
$summary.Cells.Item(4,4).Value = $summary.Cells.Item(3,4).Value2

$summary.Cells.Item(3,2).Value = $summary.Cells.Item(2,2).Value2
$summary.Cells.Item(3,3).Value = $summary.Cells.Item(2,3).Value2
$summary.Cells.Item(3,4).Value = $summary.Cells.Item(2,4).Value2

# New row 2: the 2022-Q3 summary entry.
$summary.Cells.Item(2,2).Value = "2022-Q3"
$summary.Cells.Item(2,3).Value = 13
$summary.Cells.Item(2,4).Value = 0.49

# =======================================================================
# 2. Insert a brand new worksheet "2022-Q3" right before the current
#    "2022-Q2" sheet (tab order: 总计, 2022-Q3, 2022-Q2, 2022-Q1, 2021-Q4).
# =======================================================================
$q2Sheet = $wb.Worksheets.Item(2)   # currently "2022-Q2"
$q3Sheet = $wb.Worksheets.Add($q2Sheet)
$q3Sheet.Name = "2022-Q3"

# Header row (style source: the matching header cell on the summary sheet,
# which carries the bold/centered/bordered header style).
Set-StyledValue $summary.Cells.Item(1,2) $q3Sheet.Cells.Item(1,2) "基金代码"
Set-StyledValue $summary.Cells.Item(1,2) $q3Sheet.Cells.Item(1,3) "基金名称"
Set-StyledValue $summary.Cells.Item(1,2) $q3Sheet.Cells.Item(1,4) "基金规模"
Set-StyledValue $summary.Cells.Item(1,2) $q3Sheet.Cells.Item(1,5) "股票总仓位"
Set-StyledValue $summary.Cells.Item(1,2) $q3Sheet.Cells.Item(1,6) "仓位占比"
Set-StyledValue $summary.Cells.Item(1,2) $q3Sheet.Cells.Item(1,7) "持有市值(亿元)"
Set-StyledValue $summary.Cells.Item(1,2) $q3Sheet.Cells.Item(1,8) "仓位排名"

# Data rows (fund holdings), row 2..14.
$rows = @(
    @{A=0;  B="014016"; C="中信建投品质优选一年持有期混合A"; D="4.28"; E="85.84"; F="2.44"; G="0.1044"; H=7}
    @{A=1;  B="008347"; C="中信建投价值甄选混合A";           D="3.87"; E="76.99"; F="2.26"; G="0.0875"; H=8}
    @{A=2;  B="002938"; C="中银证券健康产业灵活配置混合";     D="1.98"; E="92.72"; F="4.35"; G="0.0861"; H=8}
    @{A=3;  B="008348"; C="中信建投价值甄选混合C";           D="1.59"; E="76.99"; F="2.26"; G="0.0359"; H=8}
    @{A=4;  B="005108"; C="圆信永丰双利优选定期开放灵活配置混合"; D="0.61"; E="91.10"; F="5.79"; G="0.0353"; H=3}
    @{A=5;  B="010434"; C="红土创新医疗保健股票";             D="0.44"; E="94.49"; F="7.48"; G="0.0329"; H=3}
    @{A=6;  B="014017"; C="中信建投品质优选一年持有期混合C"; D="1.20"; E="85.84"; F="2.44"; G="0.0293"; H=7}
    @{A=7;  B="001965"; C="圆信永丰兴源灵活配置混合A";       D="0.48"; E="93.86"; F="5.76"; G="0.0276"; H=3}
    @{A=8;  B="001966"; C="圆信永丰兴源灵活配置混合C";       D="0.26"; E="93.86"; F="5.76"; G="0.0150"; H=3}
    @{A=9;  B="007468"; C="中信建投策略精选混合A";           D="0.55"; E="81.86"; F="2.51"; G="0.0138"; H=6}
    @{A=10; B="007469"; C="中信建投策略精选混合C";           D="0.38"; E="81.86"; F="2.51"; G="0.0095"; H=6}
    @{A=11; B="006274"; C="圆信永丰医药健康混合";             D="0.16"; E="93.60"; F="5.52"; G="0.0088"; H=3}
    @{A=12; B="165531"; C="信诚多策略灵活配置混合（LOF）";   D="0.69"; E="72.58"; F="1.28"; G="0.0088"; H=5}
)

$r = 2
foreach ($row in $rows) {
    Set-StyledValue $summary.Cells.Item(2,1) $q3Sheet.Cells.Item($r,1) $row.A
    Set-TextValue   $q3Sheet.Cells.Item($r,2) $row.B
    $q3Sheet.Cells.Item($r,3).Value = $row.C
    Set-TextValue   $q3Sheet.Cells.Item($r,4) $row.D
    Set-TextValue   $q3Sheet.Cells.Item($r,5) $row.E
    Set-TextValue   $q3Sheet.Cells.Item($r,6) $row.F
    Set-TextValue   $q3Sheet.Cells.Item($r,7) $row.G
    $q3Sheet.Cells.Item($r,8).Value = $row.H
    $r++
}
